# "version final de l acp"
#
# 1) Relocate the stray "_GoBack" bookmark: it currently sits just before
#    ". Les pays avec des politiques energetiques avancees sont-ils" but
#    belongs a bit further up the document, inside the first question,
#    right after "... dans l'adoption des energies renouvelables " and
#    before "(PIB, education, etc.) ?".
# 2) Justify the "Normal" style's paragraphs.
# 3) Mark a few linked character styles as Quick Styles.

$d = $word.ActiveDocument

# --- 1. Move the _GoBack bookmark -----------------------------------------

# Remove the bookmark from its current (stale) location, if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the new anchor point: right after "... renouvelables " and before
# "(PIB, education, ...".
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "expliquant les variations dans l'adoption des énergies renouvelables ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPoint = $d.Range($anchor.End, $anchor.End)
    $d.Bookmarks.Add("_GoBack", $insertPoint)
}

# --- 2. Justify the Normal style -------------------------------------------

$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Alignment = 3

# --- 3. Flag linked character styles as Quick Styles -----------------------

$quickStyleNames = @("Footer Char", "Subtitle Char", "Body Text Char")
foreach ($name in $quickStyleNames) {
    $style = $d.Styles($name)
    $style.QuickStyle = $true
}
